# Add a new program row (row 43) to the "programs" sheet describing the
# expectedPensionLetter program, based on Dolls et al. (2019).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

# Column order chosen so that new shared-string entries are appended in the
# same order as in the target workbook (A, B, F, K, I, J).
$ws.Cells.Item($row, 1).Value  = "expectedPensionLetter"
$ws.Cells.Item($row, 2).Value  = "Pension Information"
$ws.Cells.Item($row, 3).Value  = 2005
$ws.Cells.Item($row, 4).Value  = "Tax Reform"
$ws.Cells.Item($row, 5).Value  = 27
$ws.Cells.Item($row, 6).Value  = "Since 2005 the German pension administration sends out letters designed to inform about one's future expected future pension payments. These letters also highlight the link between social security contriubtions and the resulting pension entitlement. To receive such a letter, the recipient had to be at least 27 years old. This age cutoff  thus generated quasi-random variation which allows evaluating the resulting effect on earnings and retirement savings."
$ws.Cells.Item($row, 7).Value  = 2004.07
$ws.Cells.Item($row, 11).Value = "dolls2019"
$ws.Cells.Item($row, 9).Value  = "Dolls et al. (2019)"

# Link column (J) gets a real hyperlink, matching the style used by every
# other row's link cell.
$ws.Hyperlinks.Add($ws.Cells.Item($row, 10), "https://ideas.repec.org/a/eee/pubeco/v171y2019icp105-116.html") | Out-Null
$ws.Cells.Item($row, 10).Style = $ws.Cells.Item($row - 1, 10).Style

# Match the tall row height used for long-notes rows like this one.
$ws.Rows.Item($row).RowHeight = 150

# Reflect the author's final viewport/selection state.
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 6
$ws.Range("A43").Select() | Out-Null
